{"js": "// Update the 20x5 addition/subtraction worksheet table: each cell's\n// formula text is replaced with the value recorded in the commit's diff.\n// Cells are addressed by (row, col) so duplicate \"before\" values (several\n// cells originally read e.g. \"71-6=\") are each mapped to their own,\n// independent replacement instead of a blind text find/replace.\nconst newValues = [\n  [\"37+9=\", \"82+14=\", \"89-1=\", \"22+49=\", \"94-30=\"],\n  [\"9+48=\", \"34-15=\", \"42+39=\", \"20+33=\", \"48-20=\"],\n  [\"73+13=\", \"22-20=\", \"90-55=\", \"50+35=\", \"71-45=\"],\n  [\"82-10=\", \"31+58=\", \"41-36=\", \"34+50=\", \"97-77=\"],\n  [\"47-23=\", \"35+39=\", \"74-26=\", \"60-51=\", \"48+8=\"],\n  [\"79+20=\", \"3+59=\", \"12+1=\", \"61+35=\", \"94-29=\"],\n  [\"34+11=\", \"52-46=\", \"26+7=\", \"96-95=\", \"23-15=\"],\n  [\"45-31=\", \"81-48=\", \"91-66=\", \"83-18=\", \"10+14=\"],\n  [\"12+49=\", \"81-10=\", \"1+97=\", \"95-38=\", \"4+20=\"],\n  [\"40-14=\", \"57+15=\", \"57-53=\", \"9+75=\", \"52+46=\"],\n  [\"90+0=\", \"49-42=\", \"10+60=\", \"17-13=\", \"89-84=\"],\n  [\"50-13=\", \"55-50=\", \"52+21=\", \"16+60=\", \"80-4=\"],\n  [\"27+5=\", \"63-50=\", \"62-30=\", \"96+2=\", \"84-52=\"],\n  [\"48-14=\", \"88-24=\", \"99-42=\", \"0+9=\", \"70-16=\"],\n  [\"84-39=\", \"91-19=\", \"91-47=\", \"87-33=\", \"48-0=\"],\n  [\"21+45=\", \"3+73=\", \"88-52=\", \"2+2=\", \"13+22=\"],\n  [\"14-13=\", \"45+48=\", \"10+4=\", \"13+22=\", \"92-39=\"],\n  [\"25-20=\", \"27+1=\", \"19+51=\", \"33+47=\", \"38+28=\"],\n  [\"57-16=\", \"55-12=\", \"81-41=\", \"60+38=\", \"54-24=\"],\n  [\"2+96=\", \"7+75=\", \"80-46=\", \"62-41=\", \"52-5=\"],\n];\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"isNullObject,rowCount\");\nawait context.sync();\n\nif (table.isNullObject) {\n  throw new Error(\"Expected worksheet table not found in document body.\");\n}\n\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    table.getCell(r, c).value = newValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the 20x5 addition/subtraction worksheet table: each cell's\n# formula text is replaced with the value recorded in the commit's diff.\n# Cells are addressed by (row, col) 1-based indices (Word COM convention) so\n# duplicate \"before\" values (several cells originally read e.g. \"71-6=\")\n# each get their own, independent replacement instead of a blind\n# find/replace across the whole document.\n$newValues = @(\n    @(\"37+9=\", \"82+14=\", \"89-1=\", \"22+49=\", \"94-30=\"),\n    @(\"9+48=\", \"34-15=\", \"42+39=\", \"20+33=\", \"48-20=\"),\n    @(\"73+13=\", \"22-20=\", \"90-55=\", \"50+35=\", \"71-45=\"),\n    @(\"82-10=\", \"31+58=\", \"41-36=\", \"34+50=\", \"97-77=\"),\n    @(\"47-23=\", \"35+39=\", \"74-26=\", \"60-51=\", \"48+8=\"),\n    @(\"79+20=\", \"3+59=\", \"12+1=\", \"61+35=\", \"94-29=\"),\n    @(\"34+11=\", \"52-46=\", \"26+7=\", \"96-95=\", \"23-15=\"),\n    @(\"45-31=\", \"81-48=\", \"91-66=\", \"83-18=\", \"10+14=\"),\n    @(\"12+49=\", \"81-10=\", \"1+97=\", \"95-38=\", \"4+20=\"),\n    @(\"40-14=\", \"57+15=\", \"57-53=\", \"9+75=\", \"52+46=\"),\n    @(\"90+0=\", \"49-42=\", \"10+60=\", \"17-13=\", \"89-84=\"),\n    @(\"50-13=\", \"55-50=\", \"52+21=\", \"16+60=\", \"80-4=\"),\n    @(\"27+5=\", \"63-50=\", \"62-30=\", \"96+2=\", \"84-52=\"),\n    @(\"48-14=\", \"88-24=\", \"99-42=\", \"0+9=\", \"70-16=\"),\n    @(\"84-39=\", \"91-19=\", \"91-47=\", \"87-33=\", \"48-0=\"),\n    @(\"21+45=\", \"3+73=\", \"88-52=\", \"2+2=\", \"13+22=\"),\n    @(\"14-13=\", \"45+48=\", \"10+4=\", \"13+22=\", \"92-39=\"),\n    @(\"25-20=\", \"27+1=\", \"19+51=\", \"33+47=\", \"38+28=\"),\n    @(\"57-16=\", \"55-12=\", \"81-41=\", \"60+38=\", \"54-24=\"),\n    @(\"2+96=\", \"7+75=\", \"80-46=\", \"62-41=\", \"52-5=\"),\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 0; $r -lt $newValues.Count; $r++) {\n    for ($c = 0; $c -lt $newValues[$r].Count; $c++) {\n        $t.Cell($r + 1, $c + 1).Range.Text = $newValues[$r][$c]\n    }\n}\n"}
